# Generate Report for Handoff
# Adds two new tracked files to the localization-status report:
#   12b9afc9-87bd-4e53-808f-7e21aff183a9.md  (inserted alphabetically before 40f63ad1...)
#   a1d2e782-46d3-4ead-b647-d382760c2800.md  (appended at the end)
# across all three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ----------------------------------------------------------------------
# Sheet "Overview"
# ----------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")
$loO = $wsO.ListObjects.Item(1)

# Existing row 5 (previously 40f63ad1) becomes the new 12b9afc9 entry.
$wsO.Cells.Item(5,1).Value = "12b9afc9-87bd-4e53-808f-7e21aff183a9.md"
$wsO.Cells.Item(5,2).Value = "e2e\12b9afc9-87bd-4e53-808f-7e21aff183a9.md"
$wsO.Cells.Item(5,3).Value = ".md"
$wsO.Cells.Item(5,4).Value = ""
$wsO.Cells.Item(5,5).Value = "Ready for handoff"
$wsO.Cells.Item(5,6).Value = "Ready for handoff"
$wsO.Cells.Item(5,7).Value = "2016-08-20 06:48:01"
$wsO.Cells.Item(5,7).NumberFormat = $dateFmt

# Append two brand-new rows to the table (6 = 40f63ad1 moved down, 7 = a1d2e782 new).
$loO.ListRows.Add() | Out-Null
$loO.ListRows.Add() | Out-Null

$wsO.Cells.Item(6,1).Value = "40f63ad1-8856-4157-9c73-3dc753f1a36a.md"
$wsO.Cells.Item(6,2).Value = "e2e\40f63ad1-8856-4157-9c73-3dc753f1a36a.md"
$wsO.Cells.Item(6,3).Value = ".md"
$wsO.Cells.Item(6,4).Value = ""
$wsO.Cells.Item(6,5).Value = "Ready for handoff"
$wsO.Cells.Item(6,6).Value = "Ready for handoff"
$wsO.Cells.Item(6,7).Value = "2016-08-20 06:46:29"
$wsO.Cells.Item(6,7).NumberFormat = $dateFmt

$wsO.Cells.Item(7,1).Value = "a1d2e782-46d3-4ead-b647-d382760c2800.md"
$wsO.Cells.Item(7,2).Value = "e2e\a1d2e782-46d3-4ead-b647-d382760c2800.md"
$wsO.Cells.Item(7,3).Value = ".md"
$wsO.Cells.Item(7,4).Value = ""
$wsO.Cells.Item(7,5).Value = "Ready for handoff"
$wsO.Cells.Item(7,6).Value = "Ready for handoff"
$wsO.Cells.Item(7,7).Value = "2016-08-20 06:48:01"
$wsO.Cells.Item(7,7).NumberFormat = $dateFmt

# Re-style column B (hyperlink column) for rows 5-7 and rebuild hyperlinks.
$wsO.Range("B5").Style = "HyperLink"
$wsO.Range("B6").Style = "HyperLink"
$wsO.Range("B7").Style = "HyperLink"

$wsO.Hyperlinks.Delete()
$wsO.Hyperlinks.Add($wsO.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a16b9f8637001cbfcb72748d0fc83d9d2bb92b1b/e2e/c658e27f-941e-48c1-a98f-0fa0197d0362.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4b43021a9c35faf88d72d67e3f5fb55113fb03f/e2e/db443c23-8287-409c-9d20-2a37581492b6.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4b43021a9c35faf88d72d67e3f5fb55113fb03f/e2e/de5ead40-71f3-4743-945e-79657328ad01.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/17572730fe1c61d461b7677786337e4118140357/e2e/12b9afc9-87bd-4e53-808f-7e21aff183a9.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/204c89b0d2cc013f4b416dd8d009a385cc522dfb/e2e/40f63ad1-8856-4157-9c73-3dc753f1a36a.md") | Out-Null
$wsO.Hyperlinks.Add($wsO.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/17572730fe1c61d461b7677786337e4118140357/e2e/a1d2e782-46d3-4ead-b647-d382760c2800.md") | Out-Null

# ----------------------------------------------------------------------
# Helper data shared by the zh-cn / de-de sheets
# ----------------------------------------------------------------------
# Columns: A Source File Name, B File Extension, C Status, D Source Path,
#          E Priority, F Content Duplicate, G Latest Handoff File,
#          H Latest Handoff Datetime, K Latest Handback DateTime,
#          M To be localized, O Has metadata

function Fill-LangRow($ws, $row, $fileName, $xlfName, $handoffDate) {
    $ws.Cells.Item($row,1).Value = $fileName
    $ws.Cells.Item($row,2).Value = ".md"
    $ws.Cells.Item($row,3).Value = "Ready for handoff"
    $ws.Cells.Item($row,4).Value = "e2e"
    $ws.Cells.Item($row,5).Value = "ht"
    $ws.Cells.Item($row,6).Value = "False"
    $ws.Cells.Item($row,7).Value = $xlfName
    $ws.Cells.Item($row,8).Value = $handoffDate
    $ws.Cells.Item($row,8).NumberFormat = $dateFmt
    $ws.Cells.Item($row,11).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item($row,11).NumberFormat = $dateFmt
    $ws.Cells.Item($row,13).Value = "True"
    $ws.Cells.Item($row,15).Value = "False"

    $ws.Range($ws.Cells.Item($row,1), $ws.Cells.Item($row,1)).Style = "HyperLink"
    $ws.Cells.Item($row,8).NumberFormat = $dateFmt
    $ws.Cells.Item($row,11).NumberFormat = $dateFmt
}

# ----------------------------------------------------------------------
# Sheet "zh-cn"
# ----------------------------------------------------------------------
$wsZ = $wb.Worksheets.Item("zh-cn")
$loZ = $wsZ.ListObjects.Item(1)

# Insert a fresh row before row 5 (currently 40f63ad1), pushing it to row 6.
$wsZ.Rows.Item(5).Insert()
$loZ.Resize($wsZ.Range("A1:P6"))

Fill-LangRow $wsZ 5 "12b9afc9-87bd-4e53-808f-7e21aff183a9.md" "12b9afc9-87bd-4e53-808f-7e21aff183a9.fd911a55b2f681c0e4a941b46e80b5df055225e4.zh-cn.xlf" "2016-08-20 06:47:56"

# Append a1d2e782 as a new final row.
$loZ.ListRows.Add() | Out-Null
Fill-LangRow $wsZ 7 "a1d2e782-46d3-4ead-b647-d382760c2800.md" "a1d2e782-46d3-4ead-b647-d382760c2800.aac6724be4371a6c839d3df914d6b1df9a74a7c6.zh-cn.xlf" "2016-08-20 06:47:56"

$wsZ.Hyperlinks.Delete()
$wsZ.Hyperlinks.Add($wsZ.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a16b9f8637001cbfcb72748d0fc83d9d2bb92b1b/e2e/c658e27f-941e-48c1-a98f-0fa0197d0362.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e88176cc59c16dcd5a336362d8172c1da412fd47/e2e/c658e27f-941e-48c1-a98f-0fa0197d0362.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4b43021a9c35faf88d72d67e3f5fb55113fb03f/e2e/db443c23-8287-409c-9d20-2a37581492b6.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4b43021a9c35faf88d72d67e3f5fb55113fb03f/e2e/de5ead40-71f3-4743-945e-79657328ad01.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/17572730fe1c61d461b7677786337e4118140357/e2e/12b9afc9-87bd-4e53-808f-7e21aff183a9.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/204c89b0d2cc013f4b416dd8d009a385cc522dfb/e2e/40f63ad1-8856-4157-9c73-3dc753f1a36a.md") | Out-Null
$wsZ.Hyperlinks.Add($wsZ.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/17572730fe1c61d461b7677786337e4118140357/e2e/a1d2e782-46d3-4ead-b647-d382760c2800.md") | Out-Null

# ----------------------------------------------------------------------
# Sheet "de-de"
# ----------------------------------------------------------------------
$wsD = $wb.Worksheets.Item("de-de")
$loD = $wsD.ListObjects.Item(1)

$wsD.Rows.Item(5).Insert()
$loD.Resize($wsD.Range("A1:P6"))

Fill-LangRow $wsD 5 "12b9afc9-87bd-4e53-808f-7e21aff183a9.md" "12b9afc9-87bd-4e53-808f-7e21aff183a9.fd911a55b2f681c0e4a941b46e80b5df055225e4.de-de.xlf" "2016-08-20 06:48:01"

$loD.ListRows.Add() | Out-Null
Fill-LangRow $wsD 7 "a1d2e782-46d3-4ead-b647-d382760c2800.md" "a1d2e782-46d3-4ead-b647-d382760c2800.aac6724be4371a6c839d3df914d6b1df9a74a7c6.de-de.xlf" "2016-08-20 06:48:01"

$wsD.Hyperlinks.Delete()
$wsD.Hyperlinks.Add($wsD.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a16b9f8637001cbfcb72748d0fc83d9d2bb92b1b/e2e/c658e27f-941e-48c1-a98f-0fa0197d0362.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f62eb2098b6d8e6a410757dca05519b522eca19b/e2e/c658e27f-941e-48c1-a98f-0fa0197d0362.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4b43021a9c35faf88d72d67e3f5fb55113fb03f/e2e/db443c23-8287-409c-9d20-2a37581492b6.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4b43021a9c35faf88d72d67e3f5fb55113fb03f/e2e/de5ead40-71f3-4743-945e-79657328ad01.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/17572730fe1c61d461b7677786337e4118140357/e2e/12b9afc9-87bd-4e53-808f-7e21aff183a9.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/204c89b0d2cc013f4b416dd8d009a385cc522dfb/e2e/40f63ad1-8856-4157-9c73-3dc753f1a36a.md") | Out-Null
$wsD.Hyperlinks.Add($wsD.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/17572730fe1c61d461b7677786337e4118140357/e2e/a1d2e782-46d3-4ead-b647-d382760c2800.md") | Out-Null

# sheet2/sheet3 row2 H/K columns keep their original (non-string) literal
# values - column H2 (zh-cn) / K2 already correct, untouched above.
